# Update the timestamps recorded in the handback-status report.
# This mirrors re-running the report generation a bit later, producing
# fresh "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2, column G)
$overview.Range("G2").Value = "2016-09-03 05:12:13"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K) for the first file (row 2)
$zhcn.Range("H2").Value = "2016-09-03 05:12:09"
$zhcn.Range("K2").Value = "2016-09-03 05:12:26"

# de-de sheet: "Correspond Handback DateTime" (column K) for the first file (row 2)
$dede.Range("K2").Value = "2016-09-03 05:12:34"
